$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns - and the two swapped
# coin rows (Monero / PolygonEcosystemToken) - with freshly scraped values.
# A leading apostrophe forces Excel to keep purely numeric-looking price
# strings (e.g. "576.67") as text instead of auto-converting them to numbers,
# matching the inline-string text cells used throughout the sheet.

$ws.Range("D2").Value = "63.151.43"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.475.94"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'576.67"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'146.76"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "2.475.62"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'29.12"
$ws.Range("E14").Value = "  +8.70%  "
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "2.924.17"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "63.195.41"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "2.484.09"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'8.11"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "'330.66"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +9.51%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'66.46"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").Value = "'667.40"
$ws.Range("E26").Value = "  +8.87%  "
$ws.Range("D27").Value = "'9.44"
$ws.Range("E27").Value = "  +12.15%  "
$ws.Range("D28").Value = "0.0₃0996"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "2.594.92"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "'8.12"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'4.79"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'5.51"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'153.77"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.373"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "'18.79"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "0.0₆0307"
$ws.Range("E45").Value = "  +9.64%  "
$ws.Range("D46").Value = "'15.17"
$ws.Range("D47").Value = "'148.95"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "'20.92"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").Value = "'0.0516"
$ws.Range("E51").Value = "  +0.68%  "
